# Applies the cryptos.xlsx price/volume update described by the commit
# "Updated cryptos list on Sun Mar  5 05:00:04 UTC 2023 with GitHub Actions".
#
# Column D ("Price") cells are plain-text in the source workbook
# (t="inlineStr"). Several of the new Price strings parse as plain numbers
# (e.g. "16.48", "1.0000") -- if assigned directly, Excel COM would
# auto-convert them to numeric cells, losing the literal text (trailing
# zeros, exact digits). Prefixing with a single-quote forces Excel to keep
# the literal text, exactly like typing the text into a cell by hand --
# the apostrophe itself is not stored, only the text that follows it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.497.20'
$ws.Range("E2").Value = '  +0.46%  '

# Row 3
$ws.Range("D3").Value = '1.574.30'
$ws.Range("E3").Value = '  +0.08%  '

# Row 4
$ws.Range("D4").Value = '''1.0000'
$ws.Range("E4").Value = '  -0.30%  '

# Row 5
$ws.Range("E5").Value = '  -0.18%  '

# Row 6
$ws.Range("D6").Value = '''291.36'
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
$ws.Range("D7").Value = '''0.3753'
$ws.Range("E7").Value = '  -0.28%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").Value = '''0.3416'
$ws.Range("E9").Value = '  +0.17%  '

# Row 10
$ws.Range("D10").Value = '''1.152'
$ws.Range("E10").Value = '  -1.03%  '

# Row 11
$ws.Range("D11").Value = '''0.07594'
$ws.Range("E11").Value = '  -0.95%  '

# Row 12
$ws.Range("D12").Value = '''1.000'
$ws.Range("E12").Value = '  -0.34%  '

# Row 13
$ws.Range("D13").Value = '''21.45'
$ws.Range("E13").Value = '  +0.57%  '

# Row 14
$ws.Range("D14").Value = '''6.013'
$ws.Range("E14").Value = '  +0.69%  '

# Row 15
$ws.Range("D15").Value = '''6.975'
$ws.Range("E15").Value = '  +0.96%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.563.07'
$ws.Range("E16").Value = '  -0.66%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.00001127'
$ws.Range("E17").Value = '  -0.70%  '

# Row 18
$ws.Range("D18").Value = '''91.27'
$ws.Range("E18").Value = '  +0.80%  '

# Row 19
$ws.Range("D19").Value = '''0.06747'
$ws.Range("E19").Value = '  +0.22%  '

# Row 20
$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '  -0.14%  '

# Row 21
$ws.Range("D21").Value = '''6.296'
$ws.Range("E21").Value = '  +1.20%  '

# Row 22
$ws.Range("D22").Value = '''16.48'
$ws.Range("E22").Value = '  -1.42%  '

# Row 23
$ws.Range("E23").Value = '  +2.02%  '

# Row 24
$ws.Range("D24").Value = '22.476.31'
$ws.Range("E24").Value = '  +0.36%  '

# Row 25
$ws.Range("D25").Value = '''2.318'
$ws.Range("E25").Value = '  -4.29%  '

# Row 26
$ws.Range("D26").Value = '''2.619'
$ws.Range("E26").Value = '  -5.15%  '

# Row 27
$ws.Range("D27").Value = '''20.20'
$ws.Range("E27").Value = '  -0.27%  '

# Row 28
$ws.Range("D28").Value = '''148.61'
$ws.Range("E28").Value = '  +2.35%  '

# Row 29
$ws.Range("D29").Value = '''4.999'
$ws.Range("E29").Value = '  -1.00%  '

# Row 30
$ws.Range("D30").Value = '''126.23'
$ws.Range("E30").Value = '  +0.11%  '

# Row 31
$ws.Range("D31").Value = '1.740.38'
$ws.Range("E31").Value = '  -0.44%  '

# Row 32
$ws.Range("E32").Value = '  +3.49%  '

# Row 33
$ws.Range("D33").Value = '''6.183'
$ws.Range("E33").Value = '  -0.22%  '

# Row 34
$ws.Range("D34").Value = '''1.977'
$ws.Range("E34").Value = '  -2.28%  '

# Row 35
$ws.Range("D35").Value = '''9.938'
$ws.Range("E35").Value = '  -0.92%  '

# Row 36
$ws.Range("D36").Value = '''0.08463'
$ws.Range("E36").Value = '  -0.83%  '

# Row 37
$ws.Range("D37").Value = '''1.384'
$ws.Range("E37").Value = '  +2.11%  '

# Row 38
$ws.Range("D38").Value = '''0.02477'
$ws.Range("E38").Value = '  -3.25%  '

# Row 39
$ws.Range("D39").Value = '''0.2303'
$ws.Range("E39").Value = '  -0.72%  '

# Row 40
$ws.Range("D40").Value = '''0.06562'
$ws.Range("E40").Value = '  +0.66%  '

# Row 41
$ws.Range("D41").Value = '''5.518'
$ws.Range("E41").Value = '  +0.80%  '

# Row 42
$ws.Range("D42").Value = '''11.47'
$ws.Range("E42").Value = '  -1.16%  '

# Row 43
$ws.Range("D43").Value = '''0.6316'
$ws.Range("E43").Value = '  -2.51%  '

# Row 44
$ws.Range("D44").Value = '''14.11'
$ws.Range("E44").Value = '  +0.21%  '

# Row 45
$ws.Range("D45").Value = '''1.000'
$ws.Range("E45").Value = '  -0.19%  '

# Row 46
$ws.Range("D46").Value = '''3.817'
$ws.Range("E46").Value = '  +0.75%  '

# Row 47
$ws.Range("D47").Value = '''0.5906'
$ws.Range("E47").Value = '  -2.06%  '

# Row 48
$ws.Range("D48").Value = '''2.103'
$ws.Range("E48").Value = '  +0.28%  '

# Row 49
$ws.Range("D49").Value = '''130.40'
$ws.Range("E49").Value = '  +4.15%  '

# Row 50
$ws.Range("E50").Value = '  -5.67%  '

# Row 51
$ws.Range("D51").Value = '''0.07350'
$ws.Range("E51").Value = '  +0.09%  '
